# Added columns for path parameters
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Tests" (sheet1): add param:type / param:uuid columns (I, J)
# and a new "Missing Required Param" test row (row 3).
# -----------------------------------------------------------------
$tests = $wb.Worksheets.Item("Tests")

# New column widths (I, J) -> stored width 20. Excel's ColumnWidth
# property is offset from the stored OOXML width by ~0.8333 (5/6),
# so subtract that to land exactly on width="20".
$tests.Columns.Item(9).ColumnWidth = 19.166666666666668
$tests.Columns.Item(10).ColumnWidth = 19.166666666666668

# Row 1 - new headers
$tests.Range("I1").Value = "param:type"
$tests.Range("J1").Value = "param:uuid"

# Row 2 - fill in values for the new columns on the existing test
$tests.Range("I2").Value = "work"
$tests.Range("J2").Value = "12345678-1234-1234-1234-123456789abc"

# Row 3 - brand new test case
$tests.Range("A3").Value = "put-data - Missing Required Param"
$tests.Range("B3").Value = "Test PUT /data/:type/:uuid with missing required parameters"

# "true" must stay a text value (matches C2), not become a boolean.
# Prefix with an apostrophe to force text entry, then reset the
# resulting "quote prefix" style back to Normal so no stray
# formatting is left behind.
$tests.Range("C3").Value = "'true"
$tests.Range("C3").Style = "Normal"

$tests.Range("D3").Value = 400
$tests.Range("E3").Value = 10000
$tests.Range("F3").Value = 2000
$tests.Range("G3").Value = 500
$tests.Range("H3").Value = "put-data,validation"

# Empty text value (kept as an actual empty-string cell, not blank).
$tests.Range("I3").Value = "'"
$tests.Range("I3").Style = "Normal"

$tests.Range("J3").Value = "12345678-1234-1234-1234-123456789abc"

# -----------------------------------------------------------------
# Sheet "Documentation" (sheet2): insert parameter documentation
# block before the "Endpoint-Specific Notes:" section, and add a
# "Required parameters" bullet to that section.
# -----------------------------------------------------------------
$docs = $wb.Worksheets.Item("Documentation")

# Make room for the 3 new rows (17-19 become the new parameter
# descriptions block; old row 17 onward shifts down by 3).
$docs.Rows.Item(18).Insert()
$docs.Rows.Item(18).Insert()
$docs.Rows.Item(18).Insert()

$docs.Range("A17").Value = "Parameter Descriptions:"
$docs.Range("A18").Value = "param:type"
$docs.Range("B18").Value = "type parameter (string) (REQUIRED - highlighted in yellow)"
$docs.Range("A19").Value = "param:uuid"
$docs.Range("B19").Value = "uuid parameter (string) (REQUIRED - highlighted in yellow)"

# Insert one more row for the "Required parameters" bullet, right
# after the "Description:" bullet line (which is now row 24).
$docs.Rows.Item(25).Insert()
$docs.Range("A25").Value = "• Required parameters: type, uuid"
